$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A2: new phone number + left-aligned style -----------------------------
$ws.Range("A2").Value = 8667361462
$ws.Range("A2").HorizontalAlignment = -4131   # xlLeft

# --- Hyperlink on B2: re-create it so Excel records "Login@123" (the old
#     display text) in the display="" attribute, then restore the
#     Hyperlink cell style (Add() resets formatting) and finally write the
#     new password text into the cell / shared string.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:Login@123", "", "", "Login@123")
$ws.Range("B2").Style = "Hyperlink"
$ws.Range("B2").Value = "siva123"
